# Scheduled Sheets runner: refresh market-data-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) with the latest
# Universalis price pull for every Leve across all crafting-job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 316.81818
$ws.Range("I11").Value = 316.81818
$ws.Range("K11").Value = 316.81818
$ws.Range("M11").Value = -176.81818
$ws.Range("H15").Value = 599.67346
$ws.Range("I15").Value = 599.67346
$ws.Range("K15").Value = 1799.02038
$ws.Range("M15").Value = -1630.02038
$ws.Range("H51").Value = 19990
$ws.Range("I51").Value = 19990
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 19990
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = -19506
$ws.Range("M51").ClearContents()
$ws.Range("H69").Value = 6865.3076
$ws.Range("J69").Value = 6865.3076
$ws.Range("L69").Value = 20595.9228
$ws.Range("N69").Value = -22343.9228
$ws.Range("H72").Value = 6865.3076
$ws.Range("J72").Value = 6865.3076
$ws.Range("L72").Value = 61787.7684
$ws.Range("N72").Value = -70523.7684
$ws.Range("H74").Value = 7617.643
$ws.Range("I74").Value = 5499.5
$ws.Range("K74").Value = 5499.5
$ws.Range("M74").Value = -4563.5
$ws.Range("H77").Value = 7617.643
$ws.Range("I77").Value = 5499.5
$ws.Range("K77").Value = 27497.5
$ws.Range("M77").Value = -22817.5
$ws.Range("H137").Value = 3967
$ws.Range("I137").Value = 1943.5
$ws.Range("J137").Value = 6665
$ws.Range("K137").Value = 5830.5
$ws.Range("L137").Value = 19995
$ws.Range("M137").Value = -3280.5
$ws.Range("N137").Value = -25095
$ws.Range("H138").Value = 3448.8413
$ws.Range("J138").Value = 3759.4082
$ws.Range("L138").Value = 11278.2246
$ws.Range("N138").Value = -21558.2246

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1043.76
$ws.Range("I32").Value = 859.1875
$ws.Range("J32").Value = 5473.5
$ws.Range("K32").Value = 859.1875
$ws.Range("L32").Value = 5473.5
$ws.Range("M32").Value = -572.1875
$ws.Range("N32").Value = -6047.5
$ws.Range("H74").Value = 1959.6471
$ws.Range("I74").Value = 1088.4166
$ws.Range("K74").Value = 1088.4166
$ws.Range("M74").Value = -214.4166
$ws.Range("H77").Value = 1959.6471
$ws.Range("I77").Value = 1088.4166
$ws.Range("K77").Value = 5442.083000000001
$ws.Range("M77").Value = -1074.083000000001
$ws.Range("H102").Value = 3779.3076
$ws.Range("I102").Value = 3594.25
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 3594.25
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = -1972.25
$ws.Range("N102").Value = -9244
$ws.Range("H122").Value = 1366.9231
$ws.Range("I122").Value = 1324.6364
$ws.Range("K122").Value = 3973.9092
$ws.Range("M122").Value = -1523.9092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 75594
$ws.Range("I82").Value = 34500
$ws.Range("K82").Value = 34500
$ws.Range("M82").Value = -34117
$ws.Range("H85").Value = 75594
$ws.Range("I85").Value = 34500
$ws.Range("K85").Value = 34500
$ws.Range("M85").Value = -33174
$ws.Range("H86").Value = 1682.8182
$ws.Range("I86").Value = 1572.1428
$ws.Range("K86").Value = 1572.1428
$ws.Range("M86").Value = -449.1428000000001
$ws.Range("H89").Value = 1682.8182
$ws.Range("I89").Value = 1572.1428
$ws.Range("K89").Value = 7860.714
$ws.Range("M89").Value = -2244.714
$ws.Range("H99").Value = 37234.79
$ws.Range("J99").Value = 20920.166
$ws.Range("L99").Value = 20920.166
$ws.Range("N99").Value = -23916.166
$ws.Range("H105").Value = 1721.3
$ws.Range("I105").Value = 1651.625
$ws.Range("K105").Value = 1651.625
$ws.Range("M105").Value = 95.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 89996
$ws.Range("J18").Value = 89996
$ws.Range("L18").Value = 89996
$ws.Range("N18").Value = -90456
$ws.Range("H22").Value = 1256.4
$ws.Range("I22").Value = 1027.6666
$ws.Range("J22").Value = 1599.5
$ws.Range("K22").Value = 1027.6666
$ws.Range("L22").Value = 1599.5
$ws.Range("M22").Value = -677.6666
$ws.Range("N22").Value = -2299.5
$ws.Range("H31").Value = 4633.16
$ws.Range("I31").Value = 1939.8125
$ws.Range("J31").Value = 9421.333000000001
$ws.Range("K31").Value = 1939.8125
$ws.Range("L31").Value = 9421.333000000001
$ws.Range("M31").Value = -1644.8125
$ws.Range("N31").Value = -10011.333
$ws.Range("H34").Value = 4633.16
$ws.Range("I34").Value = 1939.8125
$ws.Range("J34").Value = 9421.333000000001
$ws.Range("K34").Value = 1939.8125
$ws.Range("L34").Value = 9421.333000000001
$ws.Range("M34").Value = -1737.8125
$ws.Range("N34").Value = -9825.333000000001
$ws.Range("H70").Value = 80000
$ws.Range("J70").Value = 80000
$ws.Range("L70").Value = 80000
$ws.Range("N70").Value = -80630
$ws.Range("H73").Value = 80000
$ws.Range("J73").Value = 80000
$ws.Range("L73").Value = 80000
$ws.Range("N73").Value = -82184
$ws.Range("H75").Value = 94300
$ws.Range("J75").Value = 94300
$ws.Range("L75").Value = 94300
$ws.Range("N75").Value = -96296
$ws.Range("H78").Value = 94300
$ws.Range("J78").Value = 94300
$ws.Range("L78").Value = 282900
$ws.Range("N78").Value = -292884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 92700136
$ws.Range("I4").Value = 150166700
$ws.Range("J4").Value = 6500274.5
$ws.Range("K4").Value = 450500100
$ws.Range("L4").Value = 19500823.5
$ws.Range("M4").Value = -450499988
$ws.Range("N4").Value = -19501047.5
$ws.Range("H37").Value = 97799
$ws.Range("J37").Value = 97799
$ws.Range("L37").Value = 293397
$ws.Range("N37").Value = -293621
$ws.Range("H68").Value = 10339.8
$ws.Range("I68").Value = 1234.75
$ws.Range("J68").Value = 16409.834
$ws.Range("K68").Value = 3704.25
$ws.Range("L68").Value = 49229.50199999999
$ws.Range("M68").Value = -2893.25
$ws.Range("N68").Value = -50851.50199999999
$ws.Range("H71").Value = 10339.8
$ws.Range("I71").Value = 1234.75
$ws.Range("J71").Value = 16409.834
$ws.Range("K71").Value = 11112.75
$ws.Range("L71").Value = 147688.506
$ws.Range("M71").Value = -7056.75
$ws.Range("N71").Value = -155800.506
$ws.Range("H101").Value = 6000
$ws.Range("J101").Value = 6000
$ws.Range("L101").Value = 18000
$ws.Range("N101").Value = -22868

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 50000
$ws.Range("J15").Value = 50000
$ws.Range("L15").Value = 50000
$ws.Range("N15").Value = -50576
$ws.Range("H58").Value = 32333
$ws.Range("J58").Value = 32333
$ws.Range("L58").Value = 32333
$ws.Range("N58").Value = -32887
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H102").Value = 2141.9524
$ws.Range("I102").Value = 2141.9524
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2141.9524
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = -519.9524000000001
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 572.40625
$ws.Range("I16").Value = 498.92
$ws.Range("K16").Value = 498.92
$ws.Range("M16").Value = -328.92
$ws.Range("H38").Value = 20033
$ws.Range("J38").Value = 20033
$ws.Range("L38").Value = 20033
$ws.Range("N38").Value = -20853
$ws.Range("H82").Value = 741.3333
$ws.Range("I82").Value = 416.33334
$ws.Range("J82").Value = 1066.3334
$ws.Range("K82").Value = 416.33334
$ws.Range("L82").Value = 1066.3334
$ws.Range("M82").Value = -55.33334000000002
$ws.Range("N82").Value = -1788.3334
$ws.Range("H85").Value = 741.3333
$ws.Range("I85").Value = 416.33334
$ws.Range("J85").Value = 1066.3334
$ws.Range("K85").Value = 416.33334
$ws.Range("L85").Value = 1066.3334
$ws.Range("M85").Value = 831.66666
$ws.Range("N85").Value = -3562.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 24374.334
$ws.Range("I51").Value = 14061.5
$ws.Range("K51").Value = 14061.5
$ws.Range("M51").Value = -13551.5
$ws.Range("H58").Value = 11916.667
$ws.Range("I58").Value = 11916.667
$ws.Range("K58").Value = 11916.667
$ws.Range("M58").Value = -11608.667
$ws.Range("H70").Value = 48403.332
$ws.Range("H73").Value = 48403.332
$ws.Range("H123").Value = 47214.5
$ws.Range("J123").Value = 47214.5
$ws.Range("L123").Value = 47214.5
$ws.Range("N123").Value = -57014.5
$ws.Range("H132").Value = 2056.724
$ws.Range("I132").Value = 1510.0416
$ws.Range("K132").Value = 4530.1248
$ws.Range("M132").Value = -2000.1248
